$d = $word.ActiveDocument

# --- Title: month name change (منحة المعوقين لشهر أكتوبر -> نوفمبر) ---
$d.Content.Find.Execute("أكتوبر", $false, $true, $false, $false, $false, `
                         $true, 1, $false, "نوفمبر", 2) | Out-Null

# --- Table 1 (عين تموشنت district group) ---
$t1 = $d.Tables.Item(1)

# عـيـن تـمـوشـنــت row: beneficiaries / commune amount / district total
$t1.Rows.Item(2).Cells.Item(3).Range.Text  = "894"
$t1.Rows.Item(2).Cells.Item(4).Range.Text  = "8 940 000,00"
$t1.Rows.Item(2).Cells.Item(5).Range.Text  = "10 600 000,00"

# المالح district total
$t1.Rows.Item(4).Cells.Item(5).Range.Text  = "4 950 000,00"

# شعبة اللحم row
$t1.Rows.Item(5).Cells.Item(3).Range.Text  = "154"
$t1.Rows.Item(5).Cells.Item(4).Range.Text  = "1 540 000,00"

# الـعـــامـــريــــــة row
$t1.Rows.Item(8).Cells.Item(3).Range.Text  = "194"
$t1.Rows.Item(8).Cells.Item(4).Range.Text  = "1 940 000,00"

# حاسي الفلة row
$t1.Rows.Item(9).Cells.Item(3).Range.Text  = "124"
$t1.Rows.Item(9).Cells.Item(4).Range.Text  = "1 240 000,00"

# حـمـام بـوحـجــر row (+ district total)
$t1.Rows.Item(13).Cells.Item(3).Range.Text = "360"
$t1.Rows.Item(13).Cells.Item(4).Range.Text = "3 600 000,00"
$t1.Rows.Item(13).Cells.Item(5).Range.Text = "4 760 000,00"

# عـيـن الأربـعــاء row (+ district total)
$t1.Rows.Item(17).Cells.Item(3).Range.Text = "161"
$t1.Rows.Item(17).Cells.Item(4).Range.Text = "1 610 000,00"
$t1.Rows.Item(17).Cells.Item(5).Range.Text = "3 940 000,00"

# وادي الــصـبـــاح row
$t1.Rows.Item(19).Cells.Item(3).Range.Text = "107"
$t1.Rows.Item(19).Cells.Item(4).Range.Text = "1 070 000,00"

# عـيـن الـطـلـبــــة row
$t1.Rows.Item(21).Cells.Item(3).Range.Text = "95"
$t1.Rows.Item(21).Cells.Item(4).Range.Text = "950 000,00"

$t1.Rows.Item(22).Cells.Item(3).Range.Text = "129"
$t1.Rows.Item(22).Cells.Item(4).Range.Text = "1 290 000,00"

# المجموع (table 1 grand total row; first cell spans 2 grid columns)
$t1.Rows.Item(25).Cells.Item(2).Range.Text = "3294"
$t1.Rows.Item(25).Cells.Item(3).Range.Text = "32 940 000,00"
$t1.Rows.Item(25).Cells.Item(4).Range.Text = "32 940 000,00"

# --- Table 2 (بني صاف district group) ---
$t2 = $d.Tables.Item(2)

# بـنــي صــــــــاف row (+ district total)
$t2.Rows.Item(2).Cells.Item(3).Range.Text = "474"
$t2.Rows.Item(2).Cells.Item(4).Range.Text = "4 740 000,00"
$t2.Rows.Item(2).Cells.Item(5).Range.Text = "6 130 000,00"

# ولهاصة row (+ district total)
$t2.Rows.Item(5).Cells.Item(3).Range.Text = "203"
$t2.Rows.Item(5).Cells.Item(4).Range.Text = "2 030 000,00"
$t2.Rows.Item(5).Cells.Item(5).Range.Text = "2 700 000,00"

# المجموع (table 2 total row; first cell spans 2 grid columns)
$t2.Rows.Item(7).Cells.Item(2).Range.Text = "883"
$t2.Rows.Item(7).Cells.Item(3).Range.Text = "8 830 000,00"
$t2.Rows.Item(7).Cells.Item(4).Range.Text = "8 830 000,00"

# المجموع العام (grand total row; first cell spans 2 grid columns)
$t2.Rows.Item(8).Cells.Item(2).Range.Text = "4177"
$t2.Rows.Item(8).Cells.Item(3).Range.Text = "41 770 000,00"
$t2.Rows.Item(8).Cells.Item(4).Range.Text = "41 770 000,00"

# --- Final wording: spelled-out amount (41 730 000 -> 41 770 000) ---
$d.Content.Find.Execute("واحد وأربعون مليون وسبعمئة وثلاثون ألف", $false, $true, $false, $false, $false, `
                         $true, 1, $false, "واحد وأربعون مليون وسبعمئة وسبعون ألف", 2) | Out-Null
